$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the current row 2 (existing data shifts down by 7 rows).
$ws.Rows("2:8").Insert()

# The insert copies formatting down from the header row (row 1); reset the new
# rows to the default (unstyled) look before applying the correct formats.
$ws.Range("A2:C8").ClearFormats()

# The insert also leaves behind empty B-column cells that shouldn't exist for
# these new rows (only B8 ends up with an explicit, empty, styled cell).
$ws.Range("B2:B7").ClearContents()

# Re-apply the date/time number format used throughout column A by copying it
# from an existing date cell (now row 9) so the style is reused, not duplicated.
$ws.Range("A9").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new period-band rows: dates 2024-01-23 .. 2024-01-29, phase "S".
$dates = 45314, 45315, 45316, 45317, 45318, 45319, 45320
for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = "S"
}

# Row 8 also gets an (empty) bold/centered B cell, matching the header's font
# and alignment but without its border. Build that style by copying the
# header format and then stripping the border back off.
$ws.Range("A1").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B8").Borders.LineStyle = -4142
$ws.Range("B8").Value = ""

# Update the selection to reflect the newly entered phase-status column.
$ws.Range("C2:C8").Select()
